$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Automation Tests")

# All Status values (column C, rows 2-37) are updated to "Yes"
$ws.Range("C2:C37").Value = "Yes"

# Update the view: select C2:C37 (this also clears the stale scrolled
# position that used to show row 25 at the top of the window)
[void]$ws.Range("C2:C37").Select()
